$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 26 and 27 were swapped (the row "index" in column A stays put, but
# all the match data that used to be on row 27 is now on row 26, and vice
# versa).
# ---------------------------------------------------------------------------
$ws.Range("B26").Value = 6810007
$ws.Range("F26").Value = "Eupen"
$ws.Range("G26").Value = "Club Brugge"
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = "A"
$ws.Range("K26").Value = 4.75
$ws.Range("L26").Value = 4
$ws.Range("M26").Value = 1.571
$ws.Range("N26").Value = 7
$ws.Range("O26").Value = 4.75
$ws.Range("P26").Value = 1.333
$ws.Range("Q26").Value = 1.5
$ws.Range("R26").Value = 1.875
$ws.Range("S26").Value = 1.975
$ws.Range("T26").Value = 3
$ws.Range("U26").Value = 1.925
$ws.Range("V26").Value = 1.925
$ws.Range("W26").Value = -1
$ws.Range("X26").Value = -1
$ws.Range("Y26").Value = 0.333
$ws.Range("Z26").Value = -1
$ws.Range("AA26").Value = 0.9750000000000001
$ws.Range("AB26").Value = 0.925
$ws.Range("AC26").Value = -1

$ws.Range("B27").Value = 7030334
$ws.Range("F27").Value = "Cercle Brugge"
$ws.Range("G27").Value = "Genk"
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = "A"
$ws.Range("K27").Value = 2.75
$ws.Range("L27").Value = 3.5
$ws.Range("M27").Value = 2.25
$ws.Range("N27").Value = 2.4
$ws.Range("O27").Value = 3.5
$ws.Range("P27").Value = 2.55
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = 1.85
$ws.Range("S27").Value = 2
$ws.Range("T27").Value = 3
$ws.Range("U27").Value = 1.9
$ws.Range("V27").Value = 1.95
$ws.Range("W27").Value = -1
$ws.Range("X27").Value = -1
$ws.Range("Y27").Value = 1.55
$ws.Range("Z27").Value = -1
$ws.Range("AA27").Value = 1
$ws.Range("AB27").Value = -1
$ws.Range("AC27").Value = 0.95

# ---------------------------------------------------------------------------
# Rows 164 and 165 were likewise swapped.
# ---------------------------------------------------------------------------
$ws.Range("B164").Value = 6810139
$ws.Range("F164").Value = "Anderlecht"
$ws.Range("G164").Value = "Cercle Brugge"
$ws.Range("H164").Value = 2
$ws.Range("I164").Value = 0
$ws.Range("J164").Value = "H"
$ws.Range("K164").Value = 1.869
$ws.Range("L164").Value = 3.7
$ws.Range("M164").Value = 3.8
$ws.Range("N164").Value = 2.05
$ws.Range("O164").Value = 3.6
$ws.Range("P164").Value = 3.3
$ws.Range("Q164").Value = -0.25
$ws.Range("R164").Value = 1.825
$ws.Range("S164").Value = 2.025
$ws.Range("T164").Value = 2.75
$ws.Range("U164").Value = 1.9
$ws.Range("V164").Value = 1.95
$ws.Range("W164").Value = 1.05
$ws.Range("X164").Value = -1
$ws.Range("Y164").Value = -1
$ws.Range("Z164").Value = 0.825
$ws.Range("AA164").Value = -1
$ws.Range("AB164").Value = -1
$ws.Range("AC164").Value = 0.95

$ws.Range("B165").Value = 6810143
$ws.Range("F165").Value = "Charleroi"
$ws.Range("G165").Value = "KV Mechelen"
$ws.Range("H165").Value = 3
$ws.Range("I165").Value = 1
$ws.Range("J165").Value = "H"
$ws.Range("K165").Value = 2.05
$ws.Range("L165").Value = 3.4
$ws.Range("M165").Value = 3.5
$ws.Range("N165").Value = 2.375
$ws.Range("O165").Value = 3.1
$ws.Range("P165").Value = 3.1
$ws.Range("Q165").Value = -0.25
$ws.Range("R165").Value = 2.025
$ws.Range("S165").Value = 1.825
$ws.Range("T165").Value = 2.25
$ws.Range("U165").Value = 2.025
$ws.Range("V165").Value = 1.825
$ws.Range("W165").Value = 1.375
$ws.Range("X165").Value = -1
$ws.Range("Y165").Value = -1
$ws.Range("Z165").Value = 1.025
$ws.Range("AA165").Value = -1
$ws.Range("AB165").Value = 1.025
$ws.Range("AC165").Value = -1

# ---------------------------------------------------------------------------
# Rows 245-247 get refreshed with newer data/odds (what used to live on rows
# 248-250), and the now-duplicated trailing rows 248-250 are removed.
# ---------------------------------------------------------------------------
$ws.Range("B245").Value = 7979463
$ws.Range("E245").Value = 45382.5625
$ws.Range("F245").Value = "SintTruidense"
$ws.Range("G245").Value = "Westerlo"
$ws.Range("K245").Value = 2.2
$ws.Range("L245").Value = 3.5
$ws.Range("M245").Value = 3.2
$ws.Range("N245").Value = 2.2
$ws.Range("O245").Value = 3.5
$ws.Range("P245").Value = 3.25
$ws.Range("Q245").Value = -0.25
$ws.Range("R245").Value = 1.875
$ws.Range("S245").Value = 1.975
$ws.Range("T245").Value = 2.5
$ws.Range("U245").Value = 1.875
$ws.Range("V245").Value = 1.975

$ws.Range("B246").Value = 7979348
$ws.Range("E246").Value = 45383.35416666666
$ws.Range("F246").Value = "Cercle Brugge"
$ws.Range("G246").Value = "Club Brugge"
$ws.Range("K246").Value = 3.6
$ws.Range("L246").Value = 3.6
$ws.Range("M246").Value = 2
$ws.Range("N246").Value = 3.6
$ws.Range("O246").Value = 3.6
$ws.Range("P246").Value = 2
$ws.Range("Q246").Value = 0.5
$ws.Range("R246").Value = 1.85
$ws.Range("S246").Value = 2
$ws.Range("T246").Value = 2.75
$ws.Range("U246").Value = 1.925
$ws.Range("V246").Value = 1.925

$ws.Range("B247").Value = 7979347
$ws.Range("E247").Value = 45383.5625
$ws.Range("F247").Value = "Genk"
$ws.Range("G247").Value = "Union Saint Gilloise"
$ws.Range("K247").Value = 2.8
$ws.Range("L247").Value = 3.4
$ws.Range("M247").Value = 2.45
$ws.Range("N247").Value = 2.75
$ws.Range("O247").Value = 3.4
$ws.Range("P247").Value = 2.5
$ws.Range("Q247").Value = 0
$ws.Range("R247").Value = 2.05
$ws.Range("S247").Value = 1.8
$ws.Range("T247").Value = 2.75
$ws.Range("U247").Value = 1.95
$ws.Range("V247").Value = 1.9

# Remove the rows that are no longer needed (their data now lives on
# rows 245-247 above), shrinking the table from 250 to 247 data rows.
$ws.Range("A248:AC250").Delete()
